$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '26.439.01', '  -0.17%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.811.81', '  +0.36%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.003', '  -0.63%  ')
    ,@(6, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '305.40', '  -0.99%  ')
    ,@(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4505', '  -0.76%  ')
    ,@(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3584', '  -2.01%  ')
    ,@(9, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '46.29', '  +2.69%  ')
    ,@(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07058', '  -0.72%  ')
    ,@(11, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.8891', '  +1.52%  ')
    ,@(12, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07773', '  +0.37%  ')
    ,@(13, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '19.30', '  -0.22%  ')
    ,@(14, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.788.58', '  -2.41%  ')
    ,@(15, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.269', '  +0.10%  ')
    ,@(16, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.307', '  -0.50%  ')
    ,@(17, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '84.98', '  -1.23%  ')
    ,@(18, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.005', '  -0.58%  ')
    ,@(19, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000008516', '  -0.54%  ')
    ,@(20, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.003', '  -0.64%  ')
    ,@(21, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '26.485.03', '  -0.16%  ')
    ,@(22, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '14.14', '  -0.57%  ')
    ,@(23, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.956', '  -0.26%  ')
    ,@(24, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.036.09', '  -0.12%  ')
    ,@(25, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '10.49', '  +0.92%  ')
    ,@(26, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.958', '  -0.90%  ')
    ,@(27, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '150.59', '  +0.14%  ')
    ,@(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '17.75', '  -0.84%  ')
    ,@(29, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.049', '  +2.87%  ')
    ,@(30, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '112.32', '  -0.15%  ')
    ,@(31, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.823', '  -0.17%  ')
    ,@(32, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.08676', '  +0.31%  ')
    ,@(33, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '3.138', '  +3.21%  ')
    ,@(34, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7450', '  +2.47%  ')
    ,@(35, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.735', '  +6.94%  ')
    ,@(36, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.424', '  -0.12%  ')
    ,@(37, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.107', '  -0.34%  ')
    ,@(38, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.066', '  -1.24%  ')
    ,@(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01926', '  -0.03%  ')
    ,@(40, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.900', '  +0.77%  ')
    ,@(41, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05091', '  +0.02%  ')
    ,@(42, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.5082', '  +1.68%  ')
    ,@(43, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.749', '  -2.62%  ')
    ,@(44, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1506', '  -3.75%  ')
    ,@(45, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '8.035', '  -0.87%  ')
    ,@(46, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.4698', '  +2.30%  ')
    ,@(47, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.002', '  -0.66%  ')
    ,@(48, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '9.976', '  +0.40%  ')
    ,@(49, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '100.04', '  -1.50%  ')
    ,@(50, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.573', '  -0.60%  ')
    ,@(51, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05983', '  -0.09%  ')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowData = $data[$i]
    $row = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = "'" + $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}

